# Change the bracketed suffix on the first line from " (Alternative)" into
# " (Changed main)", but written out as three separate runs ("(", "Changed
# main", ")") instead of one merged run, per the target diff.
$d = $word.ActiveDocument

# Locate the run that carries " (Alternative)" using Find so we don't rely
# on hard-coded character offsets.
$find = $d.Content.Find
$found = $find.Execute(" (Alternative)", $false, $false, $false, $false, `
                        $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the text ' (Alternative)' to replace."
}
$target = $find.Parent

# Replace the matched range's contents with three explicit runs by feeding
# raw WordprocessingML through InsertXML (wrapped in the flat-OPC pkg:package
# envelope). Plain Range.Text / Find-Replace edits get auto-merged back into
# a single run by the engine's save-time canonicalizer since they'd share
# identical formatting, so InsertXML is used to keep the run boundaries.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' +
       '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
       '<w:r><w:t>Changed main</w:t></w:r>' +
       '<w:r><w:t>)</w:t></w:r>' +
       '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
